$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10971.0735524696
$ws.Range("C2").Value = 10129.1421202102
$ws.Range("E2").Value = 6892.65258806906
$ws.Range("F2").Value = 0.230612844969376

$ws.Range("B3").Value = 10780.8225002414
$ws.Range("C3").Value = 9985.09128819418
$ws.Range("E3").Value = 6587.42121933308
$ws.Range("F3").Value = 222.510521146969

$ws.Range("B4").Value = 10723.384679257
$ws.Range("C4").Value = 9357.08566939723
$ws.Range("E4").Value = 6669.99037389533
$ws.Range("F4").Value = 199.784001803857

$ws.Range("B5").Value = 4074.00229605486
$ws.Range("C5").Value = 6638.01748304002
$ws.Range("E5").Value = 6381.51112839683
$ws.Range("F5").Value = 74.4695254765351

$ws.Range("B6").Value = 3885.53430155908
$ws.Range("C6").Value = 6484.71654042999
$ws.Range("E6").Value = 6331.31049310012
$ws.Range("F6").Value = 65.9902930637543

$ws.Range("B7").Value = 3853.32953270235
$ws.Range("C7").Value = 6398.7093312229
$ws.Range("E7").Value = 6365.96083750484
$ws.Range("F7").Value = 63.8504236969892
